# "Primer sprint de práctica 3"
# Mark the three Sprint-1 backlog items (rows 2-4, column G) as done with an
# "x", centered both horizontally and vertically (no wrap), and move the
# active selection on to the next row (B5) as the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write + format the first cell directly so only ONE new cell style
# (center/center, no wrap) gets created in cellXfs.
$g2 = $ws.Range("G2")
$g2.Value = "x"
$g2.HorizontalAlignment = -4108   # xlCenter
$g2.VerticalAlignment = -4108     # xlCenter

# Fill in the rest of the sprint-1 rows with the same mark...
$ws.Range("G3").Value = "x"
$ws.Range("G4").Value = "x"

# ...then clone G2's formatting onto them instead of re-setting the
# alignment properties (which would otherwise fork a fresh style per cell).
$g2.Copy()
$ws.Range("G3:G4").PasteSpecial(-4122)  # xlPasteFormats

# Move the selection to B5, matching where the author left off.
[void]$ws.Range("B5").Select()
